$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3845463333333334
$ws.Range("H2").Value = 1.153639
$ws.Range("I2").Value = 0.1984850200147207
$ws.Range("J2").Value = 0.1984850200147207
$ws.Range("M2").Value = 1.785326666666666
$ws.Range("N2").Value = 5.355979999999999
$ws.Range("O2").Value = 0.6170427000672803
$ws.Range("P2").Value = 0.6170427000672803
$ws.Range("Q2").Value = 0.6865408234688888
$ws.Range("R2").Value = 6.178867411219999
$ws.Range("S2").Value = 0.1224737326727915
$ws.Range("T2").Value = 0.1224737326727914
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3845463333333334
$ws.Range("H3").Value = 1.153639
$ws.Range("I3").Value = 0.1984850200147207
$ws.Range("J3").Value = 0.1984850200147207
$ws.Range("O3").Value = 0.2473424208071816
$ws.Range("P3").Value = 0.2473424208071815
$ws.Range("Q3").Value = 0.2752008398142223
$ws.Range("R3").Value = 2.476807558328
$ws.Range("S3").Value = 0.04909376534440291
$ws.Range("T3").Value = 0.0490937653444029
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3845463333333334
$ws.Range("H4").Value = 1.153639
$ws.Range("I4").Value = 0.1984850200147207
$ws.Range("J4").Value = 0.1984850200147207
$ws.Range("M4").Value = 0.3923826666666666
$ws.Range("N4").Value = 1.177148
$ws.Range("O4").Value = 0.135614879125538
$ws.Range("P4").Value = 0.135614879125538
$ws.Range("Q4").Value = 0.1508893157302222
$ws.Range("R4").Value = 1.358003841572
$ws.Range("S4").Value = 0.02691752199752634
$ws.Range("T4").Value = 0.02691752199752634
$ws.Range("I5").Value = 0.5733580031870772
$ws.Range("J5").Value = 0.5733580031870772
$ws.Range("M5").Value = 1.785326666666666
$ws.Range("N5").Value = 5.355979999999999
$ws.Range("O5").Value = 0.6170427000672803
$ws.Range("P5").Value = 0.6170427000672803
$ws.Range("Q5").Value = 1.98319085048
$ws.Range("R5").Value = 17.84871765432
$ws.Range("S5").Value = 0.3537863703917384
$ws.Range("T5").Value = 0.3537863703917384
$ws.Range("I6").Value = 0.5733580031870772
$ws.Range("J6").Value = 0.5733580031870772
$ws.Range("O6").Value = 0.2473424208071816
$ws.Range("P6").Value = 0.2473424208071815
$ws.Range("S6").Value = 0.1418157564974634
$ws.Range("T6").Value = 0.1418157564974634
$ws.Range("I7").Value = 0.5733580031870772
$ws.Range("J7").Value = 0.5733580031870772
$ws.Range("M7").Value = 0.3923826666666666
$ws.Range("N7").Value = 1.177148
$ws.Range("O7").Value = 0.135614879125538
$ws.Range("P7").Value = 0.135614879125538
$ws.Range("Q7").Value = 0.4358696528479999
$ws.Range("R7").Value = 3.922826875631999
$ws.Range("S7").Value = 0.0777558762978753
$ws.Range("T7").Value = 0.0777558762978753
$ws.Range("G8").Value = 0.442033
$ws.Range("H8").Value = 1.326099
$ws.Range("I8").Value = 0.2281569767982021
$ws.Range("J8").Value = 0.2281569767982021
$ws.Range("M8").Value = 1.785326666666666
$ws.Range("N8").Value = 5.355979999999999
$ws.Range("O8").Value = 0.6170427000672803
$ws.Range("P8").Value = 0.6170427000672803
$ws.Range("Q8").Value = 0.7891733024466664
$ws.Range("R8").Value = 7.102559722019998
$ws.Range("S8").Value = 0.1407825970027505
$ws.Range("T8").Value = 0.1407825970027504
$ws.Range("G9").Value = 0.442033
$ws.Range("H9").Value = 1.326099
$ws.Range("I9").Value = 0.2281569767982021
$ws.Range("J9").Value = 0.2281569767982021
$ws.Range("O9").Value = 0.2473424208071816
$ws.Range("P9").Value = 0.2473424208071815
$ws.Range("Q9").Value = 0.3163412111386667
$ws.Range("R9").Value = 2.847070900248
$ws.Range("S9").Value = 0.05643289896531527
$ws.Range("T9").Value = 0.05643289896531525
$ws.Range("G10").Value = 0.442033
$ws.Range("H10").Value = 1.326099
$ws.Range("I10").Value = 0.2281569767982021
$ws.Range("J10").Value = 0.2281569767982021
$ws.Range("M10").Value = 0.3923826666666666
$ws.Range("N10").Value = 1.177148
$ws.Range("O10").Value = 0.135614879125538
$ws.Range("P10").Value = 0.135614879125538
$ws.Range("Q10").Value = 0.1734460872946666
$ws.Range("R10").Value = 1.561014785652
$ws.Range("S10").Value = 0.03094148083013636
$ws.Range("T10").Value = 0.03094148083013635
